$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the review-data rows (old row 2 <-> old row 3) ---
# Column C (email), D (recovery email) and F (review text) differ between the
# two rows; columns A, B, E, G are identical across both rows already.
$ws.Range("C2").Value = "cohenyossi408@gmail.com"
$ws.Range("D2").Value = "cohenn167@gmail.com"
$ws.Range("F2").Value = "awesome app with great addictive concept"

$ws.Range("C3").Value = "cohenn167@gmail.com"
$ws.Range("D3").Value = "stavsade45@gmail.com"
$ws.Range("F3").Value = "nice car and tracks! Like it"

# --- Insert a new blank row above the old trailing formatted row (row 4),
#     pushing it down to row 5 ---
$ws.Rows(4).Insert()

# Restore row4 formatting: only C4/D4 should be blank-but-styled (matching
# the style previously carried by row 4, now shifted to row 5); clear any
# stray formatting the insert may have copied onto other cells in the row.
$ws.Rows(4).ClearFormats()
$ws.Range("C4").Style = $ws.Range("C5").Style
$ws.Range("D4").Style = $ws.Range("D5").Style

# --- Update the active selection to match the edited file (row 2 selected) ---
$ws.Rows(2).Select()
